$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "obs" column (K) values from "sig_Z_eta" to "Z_rap" for rows 2-29
$ws.Range("K2:K29").Value = "Z_rap"

# Update the active cell selection to match the new state (K32)
$ws.Range("K32").Select()
